$wb = $excel.ActiveWorkbook

# Sheets: 1 = Overview, 2 = zh-cn, 3 = de-de
$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

# 1) Update the status value "Ready for handoff" -> "In Translation".
#    This shared string is referenced by Overview!E2 & Overview!F2 (the
#    zh-cn / de-de status columns) as well as by the per-language sheets'
#    Status column (column C, row 2). Update every occurrence together so
#    the shared string gets replaced everywhere instead of forking a new one.
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# 2) Narrow the date/status columns from ~17.22 to ~13.41 characters wide.
#    Overview sheet: columns E and F.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

#    zh-cn / de-de sheets: column C.
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
